# Final update to the matrices and system methodology.
# Fixes a typo ("rmetoyer2016SME" -> "metoyer2016SME") in the @cite
# citation tags scattered across the "Heat Stroke" comparison matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Heat Stroke")

$cells = @("J3","L3","X3","L4","N4","P4","R4","V4","X4","AB4","V6","X6","V7","X7")

foreach ($cell in $cells) {
    $rng = $ws.Range($cell)
    $old = $rng.Text
    $new = $old.Replace("rmetoyer2016SME", "metoyer2016SME")
    $rng.Value = $new
}

# Leave the workbook selection where the author's last edit was made.
$ws.Range("L6").Select()
